$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Change 1: merge " from Board to the new" + " " (two separate runs) into a
# single run " from Board to the new " by doing a Find & Replace over the
# combined text (Find/Replace collapses the matched runs into one run).
# --------------------------------------------------------------------------
$d.Content.Find.Execute(
    "from Board to the new ", $true, $false, $false, $false, $false,
    $true, 1, $false, "from Board to the new ", 2) | Out-Null

# --------------------------------------------------------------------------
# Change 2: add a new bullet item after "... work with new classes",
# containing "IAcceptable interface added (Improved Visitor Design
# Pattern)", and move the _GoBack bookmark to the end of that new item.
# --------------------------------------------------------------------------

# The hidden "_GoBack" bookmark currently sits right after "work with new
# classes" (end of the last paragraph). Remove it now; it will be re-created
# at the end of the freshly inserted paragraph below.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Locate the end of "work with new classes" and insert a new paragraph right
# after it (this also duplicates the paragraph's list formatting/style).
$r = $d.Content
$r.Find.Execute(
    "work with new classes", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertParagraphAfter()

# Fill the brand-new (still empty) paragraph with the required runs,
# including the proofErr spell-check markers and the relocated bookmark.
$newPara = $d.Paragraphs.Last
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>IAcceptable</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> interface added (Improved Visitor Design Pattern)</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
       '</w:p>'
$newPara.Range.InsertXML($xml) | Out-Null
